# Update "4298091868_VCHAINS" workbook: refresh the value-chain data table,
# extend it from 5 rows to 10 rows, and update the defined names / column
# widths / headline fields to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / metadata cells -------------------------------------------------
$ws.Range("B1").Value = 44596.4854166667          # refreshed timestamp
$ws.Range("B3").Value = "YDx Innovation Corp"     # SheetTitle company name
$ws.Range("B4").Value = 4298091868                # Company Id

# --- Extend the data block: clone formatting for the five new rows ----------
# Row 7 carries the canonical per-column style (s=5 default, s=7 confidence
# score %, s=8 date, s=9 revenue accounting format) - copy that pattern down
# into rows 12-16 before writing their values.
$ws.Range("A7:N7").Copy() | Out-Null
$ws.Range("A12:N12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:N13").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:N14").PasteSpecial(-4122) | Out-Null
$ws.Range("A15:N15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:N16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 7 : Cisco Systems Inc ----------------------------------------------
$ws.Range("A7").Value = 4295905952
$ws.Range("B7").Value = "Cisco Systems Inc"
$ws.Range("C7").Value = "Public"
$ws.Range("D7").Value = "Customer"
$ws.Range("E7").Value = "United States of America"
$ws.Range("F7").Value = "Communications & Networking"
$ws.Range("G7").Value = 0.99284769890344
$ws.Range("H7").Value = 43276
$ws.Range("I7").Value = 1320
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 12
$ws.Range("L7").Value = 49818000000
$ws.Range("M7").Value = 66
$ws.Range("N7").Value = "A"

# --- Row 8 : Coca-Cola Co ----------------------------------------------------
$ws.Range("A8").Value = 4295903091
$ws.Range("B8").Value = "Coca-Cola Co"
$ws.Range("C8").Value = "Public"
$ws.Range("D8").Value = "Customer"
$ws.Range("E8").Value = "United States of America"
$ws.Range("F8").Value = "Non-Alcoholic Beverages"
$ws.Range("G8").Value = 0.985073674374549
$ws.Range("H8").Value = 43276
$ws.Range("I8").Value = 1320
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 33014000000
$ws.Range("M8").Value = 77
$ws.Range("N8").Value = "BB+"

# --- Row 9 : Qualcomm Inc ----------------------------------------------------
$ws.Range("A9").Value = 4295907706
$ws.Range("B9").Value = "Qualcomm Inc"
$ws.Range("C9").Value = "Public"
$ws.Range("D9").Value = "Customer"
$ws.Range("E9").Value = "United States of America"
$ws.Range("F9").Value = "Semiconductors"
$ws.Range("G9").Value = 0.984741442671968
$ws.Range("H9").Value = 43013
$ws.Range("I9").Value = 1583
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = 33566000000
$ws.Range("M9").Value = 68
$ws.Range("N9").Value = "BBB"

# --- Row 10 : Intel Corp -----------------------------------------------------
$ws.Range("A10").Value = 4295906830
$ws.Range("B10").Value = "Intel Corp"
$ws.Range("C10").Value = "Public"
$ws.Range("D10").Value = "Customer"
$ws.Range("E10").Value = "United States of America"
$ws.Range("F10").Value = "Semiconductors"
$ws.Range("G10").Value = 0.964516059202449
$ws.Range("H10").Value = 42997
$ws.Range("I10").Value = 1599
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 8
$ws.Range("L10").Value = 79024000000
$ws.Range("M10").Value = 49
$ws.Range("N10").Value = "BBB"

# --- Row 11 : Fortune 500 Marketing Solutions LLC ---------------------------
$ws.Range("A11").Value = 5049259010
$ws.Range("B11").Value = "Fortune 500 Marketing Solutions LLC"
$ws.Range("C11").Value = "Private"
$ws.Range("D11").Value = "Customer"
$ws.Range("E11").Value = "United States of America"
$ws.Range("F11").Value = "Advertising & Marketing"
$ws.Range("G11").Value = 0.719238101960211
$ws.Range("H11").Value = 42977
$ws.Range("I11").Value = 1619
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 7
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

# --- Row 12 : Adidas AG ------------------------------------------------------
$ws.Range("A12").Value = 4295868725
$ws.Range("B12").Value = "Adidas AG"
$ws.Range("C12").Value = "Public"
$ws.Range("D12").Value = "Customer"
$ws.Range("E12").Value = "Germany"
$ws.Range("F12").Value = "Footwear"
$ws.Range("G12").Value = 0.26988672
$ws.Range("H12").Value = 42885
$ws.Range("I12").Value = 1711
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 22651215034.0576
$ws.Range("M12").Value = 97
$ws.Range("N12").Value = "A-"

# --- Row 13 : Mercedes-Benz Bank AG -----------------------------------------
$ws.Range("A13").Value = 5000057425
$ws.Range("B13").Value = "Mercedes-Benz Bank AG"
$ws.Range("C13").Value = "Private"
$ws.Range("D13").Value = "Customer"
$ws.Range("E13").Value = "Germany"
$ws.Range("F13").Value = "Banks"
$ws.Range("G13").Value = 0.268392
$ws.Range("H13").Value = 42885
$ws.Range("I13").Value = 1711
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 584848668.800819
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# --- Row 14 : Nike Inc -------------------------------------------------------
$ws.Range("A14").Value = 4295904620
$ws.Range("B14").Value = "Nike Inc"
$ws.Range("C14").Value = "Public"
$ws.Range("D14").Value = "Customer"
$ws.Range("E14").Value = "United States of America"
$ws.Range("F14").Value = "Footwear"
$ws.Range("G14").Value = 0.26866416
$ws.Range("H14").Value = 42885
$ws.Range("I14").Value = 1711
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 44538000000
$ws.Range("M14").Value = 78
$ws.Range("N14").Value = "A-"

# --- Row 15 : Nokia Oyj -------------------------------------------------------
$ws.Range("A15").Value = 4295866480
$ws.Range("B15").Value = "Nokia Oyj"
$ws.Range("C15").Value = "Public"
$ws.Range("D15").Value = "Customer"
$ws.Range("E15").Value = "Finland"
$ws.Range("F15").Value = "Communications & Networking"
$ws.Range("G15").Value = 0.26934368
$ws.Range("H15").Value = 42885
$ws.Range("I15").Value = 1711
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 2
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 91
$ws.Range("N15").Value = "A"

# --- Row 16 : Unilever NV -----------------------------------------------------
$ws.Range("A16").Value = 4295884772
$ws.Range("B16").Value = "Unilever NV"
$ws.Range("C16").Value = "Private"
$ws.Range("D16").Value = "Customer"
$ws.Range("E16").Value = "Netherlands"
$ws.Range("F16").Value = "Personal Products"
$ws.Range("G16").Value = 0.26798336
$ws.Range("H16").Value = 42885
$ws.Range("I16").Value = 1711
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 57899628673.0264
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# --- Defined names: Data block grew from rows 7-11 to rows 7-16 -------------
$wb.Names.Item("Data").RefersTo = "='Value Chains'!`$A`$7:`$N`$16"
$wb.Names.Item("DataConfidenceScore").RefersTo = "='Value Chains'!`$G`$7:`$G`$16"
$wb.Names.Item("DataRevenue").RefersTo = "='Value Chains'!`$I`$7:`$I`$16"

# --- Column width tweaks (B, F, L widened) ----------------------------------
$ws.Columns.Item(2).ColumnWidth = 32.357291666666665
$ws.Columns.Item(6).ColumnWidth = 27.687291666666667
$ws.Columns.Item(12).ColumnWidth = 17.467291666666668
